# Commit: "Changed the script to include all the sheets, updated Readme.md"
#
# Sheet1 originally held 4 rows of data on rows 1,3,5,7 (with blank rows
# 2,4,6 between them). The edit:
#   1. Compacts Sheet1's data onto contiguous rows 1-4 (deletes the blank
#      rows), and refreshes the selection/dimension accordingly.
#   2. Adds three more sheets (Sheet2, Sheet3, Sheet4), each containing the
#      same A/C data laid out on contiguous rows 1-4.
#   3. Leaves Sheet1 as the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Compact Sheet1: remove the blank rows between the data rows ---
# Original layout: row1=data, row2=blank, row3=data, row4=blank,
#                   row5=data, row6=blank, row7=data.
# Deleting row 2 shifts everything up by one; what was row4 becomes row3,
# what was row6 becomes row5 -- repeating the delete at the next blank
# position compacts all the data onto rows 1-4.
$ws1.Rows.Item(2).Delete()
$ws1.Rows.Item(3).Delete()
$ws1.Rows.Item(4).Delete()

# Row 7's custom height ("ht=14.25") travelled down to row 4 with the
# data; AutoFit clears the explicit/custom height so it matches the
# other (default-height) rows, same as the target file.
$ws1.Rows.Item(4).AutoFit()

# Selection becomes a plain A1:C7 block select (matches the target sqref).
$ws1.Range("A1:C7").Select() | Out-Null

# Sheet1's C column gets a "best fit" (auto-sized) width, since its text
# ("wefWEFWE", ...) is wider than the default column width.
$ws1.Columns.Item(3).ColumnWidth = 10.3

# --- 2. Add Sheet2, Sheet3, Sheet4 after Sheet1, each getting a copy of
#        the same data (laid out on contiguous rows, like Sheet1 now is) ---
$data = @(
    @("watewet", "WEFWf"),
    @("WEFWFE", "wefwef"),
    @("WEFwef", "WFEWef"),
    @("WEFWEF", "wefWEFWE")
)

$prev = $ws1
for ($n = 2; $n -le 4; $n++) {
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prev)

    for ($r = 0; $r -lt $data.Length; $r++) {
        $ws.Cells.Item($r + 1, 1).Value = $data[$r][0]
        $ws.Cells.Item($r + 1, 3).Value = $data[$r][1]
    }

    $ws.Range("A1:C7").Select() | Out-Null

    $prev = $ws
}

# Re-activate Sheet1 so it remains the visible/selected tab.
$ws1.Activate()
